# Apply the roster update: replace "Tristan da Silva" (SF, Orlando Magic) with
# "P.J. Washington" (PF, Dallas Mavericks), and re-write the full player table
# (rows 2-19) to match the new row order produced by the upstream export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @("Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @("Coby White",         "PG,SG",    "Chicago Bulls"),
    @("Norman Powell",      "SG,SF",    "LA Clippers"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("P.J. Washington",    "PF",       "Dallas Mavericks"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Myles Turner",       "C",        "Indiana Pacers"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @("Cole Anthony",       "PG",       "Orlando Magic"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Walker Kessler",     "C",        "Utah Jazz"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("Jamal Murray",       "PG,SG",    "Denver Nuggets"),
    @("Brandon Ingram",     "SG,SF,PF", "New Orleans Pelicans"),
    @("Jabari Smith Jr.",   "PF,C",     "Houston Rockets"),
    @("D'Angelo Russell",   "PG",       "Brooklyn Nets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
